$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new log entry as row 66 (Sno 65), by cloning the formatting of the
# previous entry row (65) and then overwriting it with the new entry's data.
$ws.Range("A65:G65").Copy($ws.Range("A66:G66"))

$ws.Cells.Item(66, 1).Value2 = 65
$ws.Cells.Item(66, 2).Value2 = 44761
$ws.Cells.Item(66, 3).Value2 = 0.23958333333333334
$ws.Cells.Item(66, 4).Value2 = 0.29166666666666669
$ws.Cells.Item(66, 5).Formula = "=D66-C66"
$ws.Cells.Item(66, 6).Value2 = "Code"
$ws.Cells.Item(66, 7).Value2 = "1. mit b3 model 12ep train kaggle`n2. colab nb initial commit"

# The new row's description wraps onto two lines, so the row is twice as tall
# as a single-line row.
$ws.Rows.Item(66).RowHeight = 30

# Move the active selection to the next empty description cell, as in the
# authored workbook.
$ws.Range("G67").Select()
